$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1) - new column layout:
#   A = Test_Data_Type
#   B = UserName_TestData
#   C = UserName_TestData_Result
#   D = Email_TestData            (was column A)
#   E = Email_TestData_Result     (was column B)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Test_Data_Type"
$ws.Range("B1").Value = "UserName_TestData"
$ws.Range("C1").Value = "UserName_TestData_Result"
$ws.Range("D1").Value = "Email_TestData"
$ws.Range("E1").Value = "Email_TestData_Result"

# ---------------------------------------------------------------------------
# Data rows 2-9
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Both Positive"
$ws.Range("B2").Value = "Tester"
$ws.Range("C2").Value = "Testcase Passed"
$ws.Range("D2").Value = "test@test.com"
$ws.Range("E2").Value = "Testcase Passed"

$ws.Range("A3").Value = "Both blank values"
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = "Testcase Failed"
$ws.Range("E3").Value = "Testcase Failed"

$ws.Range("A4").Value = "UN blank`nEmail positive"
$ws.Range("B4").ClearContents()
$ws.Range("C4").Value = "Testcase Failed"
$ws.Range("D4").Value = "test2@test.com"
$ws.Range("E4").Value = "Testcase Passed"

$ws.Range("A5").Value = "UN positive`nEmail Blank"
$ws.Range("B5").Value = "Tester2"
$ws.Range("C5").Value = "Testcase Passed"
$ws.Range("E5").Value = "Testcase Failed"

$ws.Range("A6").Value = "UN positive`nEmail Negative"
$ws.Range("B6").Value = "Tester3"
$ws.Range("C6").Value = "Testcase Passed"
$ws.Range("D6").Value = "test@.com"
$ws.Range("E6").Value = "Testcase Failed"

$ws.Range("A7").Value = "UN positive`nEmail Negative"
$ws.Range("B7").Value = 123456
$ws.Range("C7").Value = "Testcase Passed"
$ws.Range("D7").Value = 123456
$ws.Range("E7").Value = "Testcase Failed"

$ws.Range("A8").Value = "UN positive`nEmail Negative"
$ws.Range("B8").Value = "!@#$%"
$ws.Range("C8").Value = "Testcase Passed"
$ws.Range("D8").Value = "test@1.2"
$ws.Range("E8").Value = "Testcase Passed"

$ws.Range("A9").Value = "UN positive`nEmail Negative"
$ws.Range("B9").Value = "A@#$123;"
$ws.Range("C9").Value = "Testcase Passed"
$ws.Range("D9").Value = "#$@.com"
$ws.Range("E9").Value = "Testcase Failed"

# ---------------------------------------------------------------------------
# Formatting
# ---------------------------------------------------------------------------
# Bold header row
$ws.Range("A1:E1").Font.Bold = $true

# Wrap text for the multi-line "Test_Data_Type" descriptions (A4:A9)
$ws.Range("A4:A9").WrapText = $true

# Row heights
$ws.Rows.Item(4).RowHeight = 30
$ws.Rows.Item(5).RowHeight = 37.5
$ws.Rows.Item(6).RowHeight = 30
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 30
$ws.Rows.Item(9).RowHeight = 30

# Column widths
$ws.Columns.Item(2).ColumnWidth = 23.59
$ws.Columns.Item(3).ColumnWidth = 25.75
$ws.Columns.Item(4).ColumnWidth = 27.59
$ws.Columns.Item(5).ColumnWidth = 23.25

# Selection matches the state the author left the sheet in
$ws.Range("D5").Select()

Write-Output "done"
